$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the reporting period text (shared string used by both A3 and B6)
$ws.Range("A3").Value = "1-Jul-2024 to 21-Dec-2024"
$ws.Range("B6").Value = "1-Jul-2024 to 21-Dec-2024"

# Update stock summary quantity (col B) and value (col D) figures
$ws.Range("B11").Value = 203
$ws.Range("D11").Value = 186.81
$ws.Range("B17").Value = 106.5
$ws.Range("D17").Value = 213
$ws.Range("B23").Value = 118.5
$ws.Range("D23").Value = 237
$ws.Range("B24").Value = 56.5
$ws.Range("D24").Value = 118.65
$ws.Range("B28").Value = 60.37
$ws.Range("D28").Value = 141.87
$ws.Range("B39").Value = 105.5
$ws.Range("D39").Value = 284.85000000000002
$ws.Range("B47").Value = 218
$ws.Range("D47").Value = 239.8
$ws.Range("B51").Value = 174
$ws.Range("D51").Value = 191.4
$ws.Range("B52").Value = 453
$ws.Range("D52").Value = 498.3
$ws.Range("B55").Value = 226
$ws.Range("D55").Value = 248.6
$ws.Range("B61").Value = 33.5
$ws.Range("D61").Value = 70.349999999999994
$ws.Range("B70").Value = 87
$ws.Range("D70").Value = 78.3
$ws.Range("B80").Value = 148
$ws.Range("B81").Value = 560
$ws.Range("D81").Value = 644
$ws.Range("B83").Value = 675.5
$ws.Range("D83").Value = 945.7
$ws.Range("B90").Value = 505.5
$ws.Range("D90").Value = 636.92999999999995
$ws.Range("B108").Value = 21.5
$ws.Range("D108").Value = 268.75
$ws.Range("B141").Value = 45.5
$ws.Range("D141").Value = 127.4
$ws.Range("B202").Value = 27
$ws.Range("D202").Value = 111.13
$ws.Range("B210").Value = 275.5
$ws.Range("D210").Value = 1239.75
$ws.Range("B211").Value = 128.5
$ws.Range("D211").Value = 578.25
$ws.Range("B269").Value = 70.5
$ws.Range("D269").Value = 418.77
$ws.Range("B308").Value = 16.2
$ws.Range("D308").Value = 142.56
$ws.Range("B358").Value = -10
$ws.Range("D358").Value = -100
$ws.Range("B376").Value = 25.75
$ws.Range("D376").Value = 296.13
$ws.Range("B387").Value = 1.5
$ws.Range("D387").Value = 15.68
$ws.Range("B410").Value = 31
$ws.Range("D410").Value = 86.8
$ws.Range("B413").Value = 39
$ws.Range("D413").Value = 115.05
$ws.Range("B462").Value = 28.5
$ws.Range("D462").Value = 216.6
$ws.Range("B468").Value = 18.5
$ws.Range("D468").Value = 69.38
$ws.Range("B477").Value = 1.5
$ws.Range("D477").Value = 1.65
$ws.Range("B492").Value = 87
$ws.Range("D492").Value = 198.36
$ws.Range("B524").Value = 17
$ws.Range("D524").Value = 68
$ws.Range("B536").Value = 14
$ws.Range("D536").Value = 70
$ws.Range("B537").Value = 25.5
$ws.Range("D537").Value = 127.5
$ws.Range("B546").Value = 28.5
$ws.Range("D546").Value = 163.88
$ws.Range("B561").Value = 69
$ws.Range("D561").Value = 190.44
$ws.Range("B579").Value = 246.5
$ws.Range("D579").Value = 473.4
$ws.Range("B584").Value = 32.5
$ws.Range("D584").Value = 139.1
$ws.Range("B608").Value = 40225.17
$ws.Range("D608").Value = 108750.23
